# Sock in the Dark - 5 a, b
$d = $word.ActiveDocument

# --- Hunk 1: merge the two runs in the first "Describe some test cases..."
#     paragraph (goat/cabbage problem) into a single run with trailing space.
$d.Content.Find.Execute(
    "Describe some test cases you tried out to make sure it works. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe some test cases you tried out to make sure it works. ", 2
) | Out-Null

# --- Hunk 2 & 3: remove the old "_GoBack" bookmark (around "Only if the
#     texture and length of socks are different...").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Hunk 4: after "Explain the solution in full." (Sock in the Dark
#     problem), append "  " + a colored explanation run.
$rng = $d.Content
$rng.Start = 0
$count = 0
$targetExplain = $null
while ($rng.Find.Execute("Explain the solution in full.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $count = $count + 1
    if ($count -eq 2) {
        $targetExplain = $rng.Duplicate
    }
    $rng.Collapse(0)
}

$targetExplain.Collapse(0)
$s0 = $targetExplain.Start
$insertText = "  Use your sense of touch to distinguish differences in texture and length to find a pair.  Each time you find a pair, remove it from the drawer so you can improve the odds of finding a matching pair on your next try."
$targetExplain.InsertAfter($insertText)
$e0 = $s0 + $insertText.Length

$coloredRange = $d.Range($s0 + 2, $e0)
$coloredRange.Font.Color = 3381555

# --- Hunk 5 & 6: the second "Describe some test cases..." paragraph (Sock
#     in the Dark) gets a paragraph-level color, plus "  " + a new
#     "_GoBack" bookmark + colored run with the test results.
$rng2 = $d.Content
$rng2.Start = 0
$count2 = 0
$targetDescribe = $null
while ($rng2.Find.Execute("Describe some test cases you tried out to make sure it works.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $count2 = $count2 + 1
    if ($count2 -eq 2) {
        $targetDescribe = $rng2.Duplicate
    }
    $rng2.Collapse(0)
}

$p1 = $targetDescribe.Paragraphs(1)
$paraRange = $p1.Range.Duplicate
$paraRange.MoveEnd(1, -1)
$paraRange.Text = ""

$customXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="3366FF"/></w:rPr></w:pPr><w:r><w:t>Describe some test cases you tried out to make sure it works.</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:color w:val="3366FF"/></w:rPr><w:t>While I did not have 20 pairs of the same colors, I used 20 pairs of sports socks and was able to detect differences in the feel and length to make matching pairs.  I only pulled two mismatches.</w:t></w:r></w:p></w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$paraRange.InsertXML($customXml)

# --- Hunk 7: close the bookmark right after the paragraph (as a sibling of
#     the <w:p>, not inside it) by inserting it at the paragraph boundary and
#     then removing the extra empty paragraph it creates.
$rng3 = $d.Content
$rng3.Start = 0
$rng3.Find.Execute("mismatches.") | Out-Null
$pC = $rng3.Paragraphs(1)
$boundaryPos = $pC.Range.End
$boundaryRng = $d.Range($boundaryPos, $boundaryPos)

$bmXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body><w:bookmarkEnd w:id="0"/></w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$boundaryRng.InsertXML($bmXml)

# Undo the extra empty paragraph InsertXML created for the bookmarkEnd.
$pExtra = $pC.Next()
$pExtra.Range.Delete()
